$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
